# Apply updated stats for app=Mediktor (rows 4-13, cols B:AO)
# generated from commit "updated results and code"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "B4" = 0.289
    "E4" = 0.168
    "F4" = 0.028
    "G4" = 0.168
    "H4" = 0.2
    "I4" = 0.026
    "J4" = 0.161
    "K4" = 0.307
    "L4" = 0.099
    "M4" = 0.314
    "N4" = 0.254
    "P4" = 0.149
    "Q4" = 0.478
    "R4" = 0.222
    "T4" = 0.244
    "V4" = 0.289
    "W4" = 0.242
    "X4" = 0.043
    "Z4" = 0.431
    "AA4" = 0.13
    "AB4" = 0.361
    "AC4" = 0.117
    "AE4" = 0.082
    "AF4" = 0.713
    "AH4" = 0.325
    "AI4" = 0.656
    "AJ4" = 0.165
    "AK4" = 0.406
    "AL4" = 0.671
    "AO4" = 0.68
    "B5" = 0.822
    "C5" = 0.146
    "D5" = 0.382
    "E5" = 0.667
    "F5" = 0.222
    "G5" = 0.471
    "H5" = 0.8
    "I5" = 0.16
    "J5" = 0.4
    "K5" = 0.6
    "L5" = 0.24
    "M5" = 0.49
    "N5" = 0.8
    "O5" = 0.16
    "P5" = 0.4
    "Q5" = 0.533
    "R5" = 0.249
    "S5" = 0.499
    "T5" = 0.511
    "U5" = 0.25
    "V5" = 0.5
    "W5" = 0.733
    "X5" = 0.196
    "Y5" = 0.442
    "Z5" = 0.8
    "AA5" = 0.16
    "AB5" = 0.4
    "AC5" = 0.711
    "AD5" = 0.205
    "AE5" = 0.453
    "AF5" = 0.956
    "AG5" = 0.042
    "AH5" = 0.206
    "AI5" = 0.778
    "AJ5" = 0.173
    "AK5" = 0.416
    "AL5" = 0.911
    "AM5" = 0.081
    "AN5" = 0.285
    "AO5" = 0.882
    "B6" = 0.428
    "E6" = 0.268
    "H6" = 0.32
    "K6" = 0.406
    "N6" = 0.386
    "Q6" = 0.504
    "T6" = 0.33
    "W6" = 0.364
    "Z6" = 0.5600000000000001
    "AC6" = 0.201
    "AF6" = 0.8169999999999999
    "AI6" = 0.712
    "AL6" = 0.773
    "AO6" = 0.767
    "B7" = 0.601
    "E7" = 0.418
    "H7" = 0.5
    "K7" = 0.504
    "N7" = 0.5590000000000001
    "Q7" = 0.521
    "T7" = 0.419
    "W7" = 0.521
    "Z7" = 0.6830000000000001
    "AC7" = 0.353
    "AF7" = 0.895
    "AI7" = 0.75
    "AL7" = 0.85
    "AO7" = 0.832
    "B8" = 0.742
    "C8" = 0.149
    "D8" = 0.386
    "E8" = 0.5629999999999999
    "H8" = 0.697
    "I8" = 0.158
    "J8" = 0.398
    "K8" = 0.531
    "M8" = 0.46
    "N8" = 0.713
    "O8" = 0.157
    "P8" = 0.396
    "Q8" = 0.509
    "S8" = 0.484
    "T8" = 0.445
    "W8" = 0.662
    "X8" = 0.182
    "Y8" = 0.426
    "Z8" = 0.737
    "AA8" = 0.157
    "AB8" = 0.396
    "AC8" = 0.596
    "AD8" = 0.189
    "AE8" = 0.435
    "AF8" = 0.879
    "AG8" = 0.06
    "AH8" = 0.244
    "AI8" = 0.77
    "AJ8" = 0.172
    "AK8" = 0.415
    "AL8" = 0.878
    "AM8" = 0.08599999999999999
    "AN8" = 0.294
    "AO8" = 0.842
    "B9" = 0.644
    "C9" = 0.229
    "D9" = 0.479
    "E9" = 0.444
    "F9" = 0.247
    "G9" = 0.497
    "H9" = 0.578
    "I9" = 0.244
    "J9" = 0.494
    "K9" = 0.444
    "L9" = 0.247
    "M9" = 0.497
    "N9" = 0.6
    "O9" = 0.24
    "P9" = 0.49
    "Q9" = 0.467
    "T9" = 0.356
    "U9" = 0.229
    "V9" = 0.479
    "W9" = 0.556
    "X9" = 0.247
    "Y9" = 0.497
    "Z9" = 0.644
    "AA9" = 0.229
    "AB9" = 0.479
    "AC9" = 0.489
    "AD9" = 0.25
    "AE9" = 0.5
    "AF9" = 0.756
    "AG9" = 0.185
    "AH9" = 0.43
    "AI9" = 0.756
    "AJ9" = 0.185
    "AK9" = 0.43
    "AL9" = 0.822
    "AM9" = 0.146
    "AN9" = 0.382
    "AO9" = 0.778
    "B10" = 0.778
    "C10" = 0.173
    "D10" = 0.416
    "E10" = 0.6
    "F10" = 0.24
    "G10" = 0.49
    "H10" = 0.733
    "I10" = 0.196
    "J10" = 0.442
    "K10" = 0.6
    "L10" = 0.24
    "M10" = 0.49
    "N10" = 0.778
    "O10" = 0.173
    "P10" = 0.416
    "Q10" = 0.533
    "R10" = 0.249
    "S10" = 0.499
    "T10" = 0.511
    "U10" = 0.25
    "V10" = 0.5
    "W10" = 0.733
    "X10" = 0.196
    "Y10" = 0.442
    "Z10" = 0.8
    "AA10" = 0.16
    "AB10" = 0.4
    "AC10" = 0.6
    "AD10" = 0.24
    "AE10" = 0.49
    "AF10" = 0.956
    "AG10" = 0.042
    "AH10" = 0.206
    "AI10" = 0.778
    "AJ10" = 0.173
    "AK10" = 0.416
    "AL10" = 0.911
    "AM10" = 0.081
    "AN10" = 0.285
    "AO10" = 0.882
    "B11" = 0.822
    "C11" = 0.146
    "D11" = 0.382
    "E11" = 0.667
    "F11" = 0.222
    "G11" = 0.471
    "H11" = 0.8
    "I11" = 0.16
    "J11" = 0.4
    "K11" = 0.6
    "L11" = 0.24
    "M11" = 0.49
    "N11" = 0.8
    "O11" = 0.16
    "P11" = 0.4
    "Q11" = 0.533
    "R11" = 0.249
    "S11" = 0.499
    "T11" = 0.511
    "U11" = 0.25
    "V11" = 0.5
    "W11" = 0.733
    "X11" = 0.196
    "Y11" = 0.442
    "Z11" = 0.8
    "AA11" = 0.16
    "AB11" = 0.4
    "AC11" = 0.644
    "AD11" = 0.229
    "AE11" = 0.479
    "AF11" = 0.956
    "AG11" = 0.042
    "AH11" = 0.206
    "AI11" = 0.778
    "AJ11" = 0.173
    "AK11" = 0.416
    "AL11" = 0.911
    "AM11" = 0.081
    "AN11" = 0.285
    "AO11" = 0.882
    "B12" = 1.378
    "C12" = 0.668
    "D12" = 0.8169999999999999
    "E12" = 1.633
    "F12" = 1.032
    "G12" = 1.016
    "H12" = 1.556
    "I12" = 1.191
    "J12" = 1.091
    "K12" = 1.407
    "L12" = 0.538
    "M12" = 0.733
    "N12" = 1.389
    "O12" = 0.571
    "P12" = 0.756
    "Z12" = 1.25
    "AA12" = 0.299
    "AB12" = 0.546
    "AC12" = 2
    "AD12" = 3.812
    "AE12" = 1.953
    "AF12" = 1.233
    "AG12" = 0.225
    "AH12" = 0.474
    "AJ12" = 0.028
    "AK12" = 0.167
    "AL12" = 1.098
    "AM12" = 0.08799999999999999
    "AN12" = 0.297
    "AO12" = 1.12
    "B13" = 3.533
    "C13" = 1.404
    "D13" = 1.185
    "E13" = 4.564
    "F13" = 0.707
    "G13" = 0.841
    "H13" = 4.524
    "I13" = 0.916
    "J13" = 0.957
    "K13" = 2.3
    "L13" = 0.61
    "M13" = 0.781
    "N13" = 3.333
    "O13" = 0.756
    "P13" = 0.869
    "Z13" = 2.833
    "AA13" = 3.901
    "AB13" = 1.975
    "AC13" = 6.273
    "AD13" = 2.88
    "AE13" = 1.697
    "AF13" = 1.667
    "AG13" = 0.8
    "AH13" = 0.894
    "AI13" = 1.311
    "AJ13" = 0.348
    "AK13" = 0.59
    "AL13" = 1.689
    "AM13" = 0.792
    "AN13" = 0.89
    "AO13" = 1.556
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output ("Updated {0} cells" -f $updates.Count)
